$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.669.50'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '3.440.20'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.92%  '
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("E11").Value = '  +3.06%  '
$ws.Range("D12").Value = '4.031.35'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("E13").Value = '  +2.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.25%  '
$ws.Range("D15").Value = '3.434.94'
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '62.705.71'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("E18").Value = '  +0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '386.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").Value = '3.581.57'
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  -5.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.24'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.10%  '
$ws.Range("E35").Value = '  +3.35%  '
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '31.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0772'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '2.562.99'
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("E48").Value = '  +2.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.56%  '
$ws.Range("E51").Value = '  -0.04%  '
